$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 35) to the supplier table, mirroring the
# existing rows' layout. Values are forced to text via a leading
# apostrophe (quote-prefix) so Excel doesn't auto-coerce numeric-looking
# strings (e.g. "10") or otherwise reinterpret the content, and the
# style is reset back to "Normal" right after so no stray number-format /
# quote-prefix style sticks to the cell.
function Set-TextCell {
    param($row, $col, [string]$text)
    $ws.Cells.Item($row, $col).Value = "'" + $text
    $ws.Cells.Item($row, $col).Style = "Normal"
}

$newRow = 35

Set-TextCell $newRow 1  "ass"
Set-TextCell $newRow 2  "ass"
Set-TextCell $newRow 3  ""
Set-TextCell $newRow 4  ""
Set-TextCell $newRow 5  ""
Set-TextCell $newRow 6  "ass"
Set-TextCell $newRow 7  ""
Set-TextCell $newRow 8  ""
Set-TextCell $newRow 9  ""
Set-TextCell $newRow 10 "2000 - Caisse - "
Set-TextCell $newRow 11 "1402 - Intitulé du compte - 14"
Set-TextCell $newRow 12 "0.0 (%)"
Set-TextCell $newRow 13 ""
Set-TextCell $newRow 14 "10"
Set-TextCell $newRow 15 ""
Set-TextCell $newRow 16 ""
Set-TextCell $newRow 17 ""
Set-TextCell $newRow 18 ""
Set-TextCell $newRow 19 ""
